# 2022-02-25 Borden Tool update
#
# - Concentration_Time_Data: MW-9 (column M) concentration series was
#   re-keyed so it mirrors the MW-1 (column E) series for every sampling
#   event row (rows 3-22). Rows that previously had a stray value now sit
#   blank, and rows that previously were blank now carry the MW-1 value.
# - Monitoring_Well_Information: the Latitude/Longitude of every well
#   (MW-1 .. MW-9, rows 2-10) was corrected to a tighter cluster of
#   coordinates.
# - Selections on the two data sheets were left where the author's last
#   click landed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Concentration_Time_Data  (MW-9 / column M now mirrors MW-1 / column E)
# ---------------------------------------------------------------------
$wsConc = $wb.Worksheets.Item("Concentration_Time_Data")

$wsConc.Range("M3").Value = 37.1
$wsConc.Range("M4").Value = 41.9
$wsConc.Range("M5").Value = 13
$wsConc.Range("M6").Value = 5.1
$wsConc.Range("M7").Value = 11.5
$wsConc.Range("M8").Value = 5
$wsConc.Range("M9").ClearContents()
$wsConc.Range("M10").Value = 4.6
$wsConc.Range("M11").Value = 1.85
$wsConc.Range("M12").ClearContents()
$wsConc.Range("M13").Value = 1.8
$wsConc.Range("M14").Value = 1.2
$wsConc.Range("M15").Value = 1
$wsConc.Range("M16").Value = 1.2
$wsConc.Range("M17").Value = 1
$wsConc.Range("M18").Value = 0.8
$wsConc.Range("M20").Value = 0.7
$wsConc.Range("M21").Value = 0.5

# Leave the selection where the author's session ended, without leaving
# this sheet active (Monitoring_Well_Information stays the active tab).
$wsConc.Range("I3:I22").Select()

# ---------------------------------------------------------------------
# Monitoring_Well_Information  (corrected Lat/Long per well)
# ---------------------------------------------------------------------
$wsWells = $wb.Worksheets.Item("Monitoring_Well_Information")

$wsWells.Range("B2").Value = 29.73166
$wsWells.Range("C2").Value = -95.4126

$wsWells.Range("B3").Value = 29.73297
$wsWells.Range("C3").Value = -95.41398

$wsWells.Range("B4").Value = 29.73308
$wsWells.Range("C4").Value = -95.41302

$wsWells.Range("B5").Value = 29.73268
$wsWells.Range("C5").Value = -95.41181

$wsWells.Range("B6").Value = 29.73373
$wsWells.Range("C6").Value = -95.41224

$wsWells.Range("B7").Value = 29.73325
$wsWells.Range("C7").Value = -95.41139

$wsWells.Range("B8").Value = 29.73399
$wsWells.Range("C8").Value = -95.41093

$wsWells.Range("B9").Value = 29.73269
$wsWells.Range("C9").Value = -95.41318

$wsWells.Range("B10").Value = 29.73256
$wsWells.Range("C10").Value = -95.41235

# Monitoring_Well_Information was (and remains) the active tab, so it is
# activated last and its selection is applied last too.
$wsWells.Activate()
$wsWells.Range("E15").Select()
